$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row above the current row 586, shifting the
# existing rows 586:635 down to 587:636 (dimension grows to R636).
$ws.Rows.Item(586).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A586").Value = 9
$ws.Range("B586").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C586").Value = "Metropolitana"
$ws.Range("D586").Value = 45223
$ws.Range("E586").Value = 13
$ws.Range("F586").Value = 100112044
$ws.Range("G586").Value = "Perejil"
$ws.Range("H586").Value = "Sin especificar"
$ws.Range("I586").Value = "Primera"
$ws.Range("J586").Value = 70
$ws.Range("K586").Value = 15000
$ws.Range("L586").Value = 16000
$ws.Range("M586").Value = 15500
$ws.Range("N586").Value = "`$/docena de atados"
$ws.Range("O586").Value = "Región Metropolitana"
$ws.Range("P586").Value = 5167
$ws.Range("Q586").Value = 3
$ws.Range("R586").Value = "Hortaliza"
